# Update "Pais" (countries) COVID-19 stats sheet with a newer data pull.
# The sheet is sorted descending by column B (Casos totales). A handful of
# countries' totals changed enough to jump over their neighbours, so those
# rows are rewritten (label + values) in their new sorted position while
# the row that got overtaken keeps the displaced country's old figures.
# Everything else is a straight value refresh in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 09:17"

# India (row 6) - value refresh, no reordering
$ws.Range("B6").Value = 769150
$ws.Range("C6").Value = 98
$ws.Range("D6").Value = 476565
$ws.Range("E6").Value = 271434
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 21151

# Alemania (row 19) - value refresh, no reordering
$ws.Range("D19").Value = 183600
$ws.Range("E19").Value = 6050

# Ucrania overtakes Paises Bajos: Ucrania moves up to row 37 (new data),
# Paises Bajos drops to row 38 keeping its old (unchanged) figures.
$ws.Range("A37").Value = "Ucrania"
$ws.Range("B37").Value = 51224
$ws.Range("C37").Value = 810
$ws.Range("D37").Value = 23784
$ws.Range("E37").Value = 26113
$ws.Range("G37").Value = 21
$ws.Range("H37").Value = 1327

$ws.Range("A38").Value = "Paises Bajos"
$ws.Range("B38").Value = 50746
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("H38").Value = 6135

# Armenia overtakes both Nigeria and Rumania: Armenia moves up to row 51
# (new data); Nigeria and Rumania each drop one row, keeping their old
# (unchanged) figures.
$ws.Range("A51").Value = "Armenia"
$ws.Range("B51").Value = 30346
$ws.Range("C51").Value = 526
$ws.Range("D51").Value = 18000
$ws.Range("E51").Value = 11811
$ws.Range("G51").Value = 14
$ws.Range("H51").Value = 535

$ws.Range("A52").Value = "Nigeria"
$ws.Range("B52").Value = 30249
$ws.Range("D52").Value = 12373
$ws.Range("E52").Value = 17192
$ws.Range("H52").Value = 684

$ws.Range("A53").Value = "Rumania"
$ws.Range("B53").Value = 30175
$ws.Range("D53").Value = 20799
$ws.Range("E53").Value = 7559
$ws.Range("H53").Value = 1817

# Uzbekistan (row 71) - value refresh, no reordering
$ws.Range("B71").Value = 11259
$ws.Range("C71").Value = 167
$ws.Range("E71").Value = 4152
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 47

# Australia overtakes Noruega: Australia moves up to row 73 (new data),
# Noruega drops to row 74 keeping its old (unchanged) figures.
$ws.Range("A73").Value = "Australia"
$ws.Range("B73").Value = 9059
$ws.Range("C73").Value = 173
$ws.Range("D73").Value = 7573
$ws.Range("E73").Value = 1380
$ws.Range("H73").Value = 106

$ws.Range("A74").Value = "Noruega"
$ws.Range("B74").Value = 8950
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 8138
$ws.Range("E74").Value = 561
$ws.Range("H74").Value = 251

# El Salvador (row 77) - value refresh, no reordering
$ws.Range("D77").Value = 5160
$ws.Range("E77").Value = 3163
$ws.Range("G77").Value = 8
$ws.Range("H77").Value = 243

# Tailandia (row 102) - value refresh, no reordering
$ws.Range("B102").Value = 3202
$ws.Range("C102").Value = 5
$ws.Range("D102").Value = 3085
$ws.Range("E102").Value = 59

# Georgia (row 142) - value refresh, no reordering
$ws.Range("B142").Value = 968
$ws.Range("C142").Value = 5
$ws.Range("D142").Value = 844
$ws.Range("E142").Value = 109

# Seychelles/Lesoto are tied (identical figures) - swap which label sits
# on which row to match the refreshed ranking.
$ws.Range("A184").Value = "Lesoto"
$ws.Range("A185").Value = "Seychelles"

# Gambia overtakes Polinesia Francesa: Gambia moves up to row 190 (new
# data), Polinesia Francesa drops to row 191 keeping its old figures.
$ws.Range("A190").Value = "Gambia"
$ws.Range("B190").Value = 63
$ws.Range("C190").Value = 2
$ws.Range("D190").Value = 32
$ws.Range("E190").Value = 28
$ws.Range("H190").Value = 3

$ws.Range("A191").Value = "Polinesia Francesa"
$ws.Range("B191").Value = 62
$ws.Range("D191").Value = 60
$ws.Range("E191").Value = 2
$ws.Range("H191").Value = 0

# Islas Malvinas/Groenlandia are tied (identical figures) - swap which
# label sits on which row to match the refreshed ranking.
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
